# The Pearson/BTec logo pictures in the footers/header had their
# internal "name" label swapped (image2.png <-> image1.png in the two
# footers that hold the Pearson logo, image1.jpg -> image2.jpg in the
# header that holds the BTec logo). The picture bytes / relationships
# themselves are untouched - only the cosmetic shape name changes.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineLogo($range, $newName) {
    $shape = $range.InlineShapes(1)
    # Route the rename through Selection - renaming the InlineShape
    # directly is unreliable for footer stories in this host, but goes
    # through fine once the shape is the active selection.
    [void]$shape.Select()
    $word.Selection.InlineShapes(1).Name = $newName
}

# Footer 1 (first page footer): Pearson logo, image2.png -> image1.png
Rename-InlineLogo $sec.Footers(1).Range "image1.png"

# Footer 2 (other pages footer): Pearson logo, image2.png -> image1.png
Rename-InlineLogo $sec.Footers(2).Range "image1.png"

# Header 2 (other pages header): BTec logo, image1.jpg -> image2.jpg
Rename-InlineLogo $sec.Headers(2).Range "image2.jpg"

Write-Host "Renamed 3 inline logo shapes"
